$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.043451070785522
$ws.Range("B1").Value = 1.724144220352173
$ws.Range("C1").Value = 5.151397705078125
$ws.Range("D1").Value = 1.249244928359985
$ws.Range("E1").Value = 0.3392289876937866
